$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price/volume refresh + rank-39/40 swap)
$ws.Range("D2").Value = "70.479.26"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "2.518.68"
$ws.Range("E3").Value = "  -5.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'576.24"
$ws.Range("E5").Value = "  -3.44%  "
$ws.Range("D6").Value = "'169.43"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "2.520.49"
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").Value = "'4.84"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "2.981.74"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").Value = "70.343.82"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "'25.02"
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").Value = "2.549.05"
$ws.Range("E18").Value = "  -3.89%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  -7.54%  "
$ws.Range("D20").Value = "'7.66"
$ws.Range("E20").Value = "  +6.51%  "
$ws.Range("D21").Value = "'359.60"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -5.67%  "
$ws.Range("D23").Value = "'1.97"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'69.24"
$ws.Range("E25").Value = "  -3.76%  "
$ws.Range("E26").Value = "  -5.86%  "
$ws.Range("D27").Value = "'9.14"
$ws.Range("E27").Value = "  -7.01%  "
$ws.Range("D28").Value = "2.650.79"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "0.0₃0914"
$ws.Range("E30").Value = "  -6.02%  "
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "'483.57"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'157.03"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'18.90"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'18.62"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'4.75"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.321"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").Value = "'1.65"
$ws.Range("E43").Value = "  -6.70%  "
$ws.Range("E44").Value = "  -13.04%  "
$ws.Range("D45").Value = "'2.39"
$ws.Range("E45").Value = "  -8.05%  "
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'143.37"
$ws.Range("E47").Value = "  -8.10%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("E50").Value = "  -6.36%  "
$ws.Range("E51").Value = "  -1.15%  "
